# Fix header labels on existing sheets to match the updated PO-forecast naming scheme
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend"
$newSheet = $wb.Worksheets.Add($null, $wsMonthly)
$newSheet.Name = "PO Forecast"

# Match the page margins used on the other sheets (0.75/0.75/1/1/0.5/0.5 in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Headers
$newSheet.Cells.Item(1,1).Value = "ds"
$newSheet.Cells.Item(1,2).Value = "PO_Forecast"
$newSheet.Cells.Item(1,3).Value = "yhat_lower"
$newSheet.Cells.Item(1,4).Value = "yhat_upper"

$newSheet.Cells.Item(2,1).Value = 44934.99999999999
$newSheet.Cells.Item(2,2).Value = 150
$newSheet.Cells.Item(2,3).Value = -40.13017343155859
$newSheet.Cells.Item(2,4).Value = 356.7581331664653
$newSheet.Cells.Item(3,1).Value = 44941.99999999999
$newSheet.Cells.Item(3,2).Value = 151
$newSheet.Cells.Item(3,3).Value = -57.60254789211113
$newSheet.Cells.Item(3,4).Value = 333.3345036022259
$newSheet.Cells.Item(4,1).Value = 44948.99999999999
$newSheet.Cells.Item(4,2).Value = 152
$newSheet.Cells.Item(4,3).Value = -49.78511167804821
$newSheet.Cells.Item(4,4).Value = 352.01648489617
$newSheet.Cells.Item(5,1).Value = 44997.99999999999
$newSheet.Cells.Item(5,2).Value = 160
$newSheet.Cells.Item(5,3).Value = -37.32053147012144
$newSheet.Cells.Item(5,4).Value = 357.5935357067503
$newSheet.Cells.Item(6,1).Value = 45004.99999999999
$newSheet.Cells.Item(6,2).Value = 161
$newSheet.Cells.Item(6,3).Value = -44.35889274368235
$newSheet.Cells.Item(6,4).Value = 360.1267615459504
$newSheet.Cells.Item(7,1).Value = 45018.99999999999
$newSheet.Cells.Item(7,2).Value = 163
$newSheet.Cells.Item(7,3).Value = -23.94374290368384
$newSheet.Cells.Item(7,4).Value = 347.8754159968612
$newSheet.Cells.Item(8,1).Value = 45025.99999999999
$newSheet.Cells.Item(8,2).Value = 164
$newSheet.Cells.Item(8,3).Value = -35.69318393163103
$newSheet.Cells.Item(8,4).Value = 355.885578232276
$newSheet.Cells.Item(9,1).Value = 45074.99999999999
$newSheet.Cells.Item(9,2).Value = 172
$newSheet.Cells.Item(9,3).Value = -38.66546971511478
$newSheet.Cells.Item(9,4).Value = 369.3778588907911
$newSheet.Cells.Item(10,1).Value = 45088.99999999999
$newSheet.Cells.Item(10,2).Value = 174
$newSheet.Cells.Item(10,3).Value = -23.31039426451597
$newSheet.Cells.Item(10,4).Value = 377.3962650432821
$newSheet.Cells.Item(11,1).Value = 45095.99999999999
$newSheet.Cells.Item(11,2).Value = 175
$newSheet.Cells.Item(11,3).Value = -24.87094073073287
$newSheet.Cells.Item(11,4).Value = 383.0919866129562
$newSheet.Cells.Item(12,1).Value = 45116.99999999999
$newSheet.Cells.Item(12,2).Value = 178
$newSheet.Cells.Item(12,3).Value = -19.37450378810283
$newSheet.Cells.Item(12,4).Value = 367.511132871986
$newSheet.Cells.Item(13,1).Value = 45123.99999999999
$newSheet.Cells.Item(13,2).Value = 180
$newSheet.Cells.Item(13,3).Value = -19.5223298883834
$newSheet.Cells.Item(13,4).Value = 368.0955209481625
$newSheet.Cells.Item(14,1).Value = 45130.99999999999
$newSheet.Cells.Item(14,2).Value = 181
$newSheet.Cells.Item(14,3).Value = -9.323059173376233
$newSheet.Cells.Item(14,4).Value = 374.0005200373816
$newSheet.Cells.Item(15,1).Value = 45137.99999999999
$newSheet.Cells.Item(15,2).Value = 182
$newSheet.Cells.Item(15,3).Value = -13.44797349709881
$newSheet.Cells.Item(15,4).Value = 381.914150547655
$newSheet.Cells.Item(16,1).Value = 45144.99999999999
$newSheet.Cells.Item(16,2).Value = 183
$newSheet.Cells.Item(16,3).Value = -17.75711339796704
$newSheet.Cells.Item(16,4).Value = 378.7350152790266
$newSheet.Cells.Item(17,1).Value = 45151.99999999999
$newSheet.Cells.Item(17,2).Value = 184
$newSheet.Cells.Item(17,3).Value = -13.87133666815923
$newSheet.Cells.Item(17,4).Value = 381.3758696716991
$newSheet.Cells.Item(18,1).Value = 45158.99999999999
$newSheet.Cells.Item(18,2).Value = 185
$newSheet.Cells.Item(18,3).Value = -13.38749042037016
$newSheet.Cells.Item(18,4).Value = 392.6269090153616
$newSheet.Cells.Item(19,1).Value = 45165.99999999999
$newSheet.Cells.Item(19,2).Value = 186
$newSheet.Cells.Item(19,3).Value = -4.318393784740031
$newSheet.Cells.Item(19,4).Value = 372.2580871093069
$newSheet.Cells.Item(20,1).Value = 45172.99999999999
$newSheet.Cells.Item(20,2).Value = 187
$newSheet.Cells.Item(20,3).Value = -2.366851221371205
$newSheet.Cells.Item(20,4).Value = 385.3324444805492
$newSheet.Cells.Item(21,1).Value = 45179.99999999999
$newSheet.Cells.Item(21,2).Value = 188
$newSheet.Cells.Item(21,3).Value = -12.81772046307768
$newSheet.Cells.Item(21,4).Value = 385.9291909614568

# Match the header / date formatting used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A21").PasteSpecial(-4122)

[void]$wsWeekly.Select()
[void]$wsWeekly.Range("A1").Select()
